$d = $word.ActiveDocument

# 1. "...containing every city and every road." -> "...one containing every city and one containing every road."
$d.Content.Find.Execute(
    "two text files, containing every city and every road.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "two text files, one containing every city and one containing every road.",
    2) | Out-Null

# 2. "...there are exactly two cities input, the two cities given..." -> "...and the two cities given..."
$d.Content.Find.Execute(
    "there are exactly two cities input, the two cities given",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "there are exactly two cities input, and the two cities given",
    2) | Out-Null

# 3. "...impossible for there to be any edges with a negative weight." -> "...and is therefore a good choice."
$d.Content.Find.Execute(
    "impossible for there to be any edges with a negative weight.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "impossible for there to be any edges with a negative weight and is therefore a good choice.",
    2) | Out-Null

# 4. Remove yellow highlighting from the two paragraphs that previously had it
#    (but leave the following empty paragraph's highlight untouched).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt.StartsWith("The graph class was implemented using a vector") -or
        $txt.StartsWith("The only problem we ran into during implementation")) {
        $p.Range.Font.HighlightColorIndex = 0
    }
}
